$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K6").Value = 1197.54546
$ws.Range("M6").Value = -1085.54546
$ws.Range("I6").Value = 399.18182
$ws.Range("H6").Value = 399.18182
$ws.Range("J64").Value = 9405.571
$ws.Range("N64").Value = -9901.571
$ws.Range("K64").Value = 5301.4
$ws.Range("L64").Value = 9405.571
$ws.Range("I64").Value = 5301.4
$ws.Range("M64").Value = -5053.4
$ws.Range("H64").Value = 7695.5
$ws.Range("J67").Value = 9405.571
$ws.Range("I67").Value = 5301.4
$ws.Range("N67").Value = -11121.571
$ws.Range("L67").Value = 9405.571
$ws.Range("M67").Value = -4443.4
$ws.Range("K67").Value = 5301.4
$ws.Range("H67").Value = 7695.5
$ws.Range("N86").Value = -7867.3335
$ws.Range("L86").Value = 5621.3335
$ws.Range("J86").Value = 5621.3335
$ws.Range("H86").Value = 5136.4443
$ws.Range("H89").Value = 5136.4443
$ws.Range("J89").Value = 5621.3335
$ws.Range("L89").Value = 28106.6675
$ws.Range("N89").Value = -39338.6675
$ws.Range("H107").Value = 1502.7858
$ws.Range("K107").Value = 1760.5555
$ws.Range("I107").Value = 1760.5555
$ws.Range("M107").Value = 159.4445000000001
$ws.Range("J130").Value = 19997.143
$ws.Range("L130").Value = 19997.143
$ws.Range("N130").Value = -30037.143
$ws.Range("H130").Value = 19997.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N21").Value = -8145.25
$ws.Range("L21").Value = 7397.25
$ws.Range("I21").Value = 4781.5
$ws.Range("K21").Value = 4781.5
$ws.Range("H21").Value = 5528.857
$ws.Range("M21").Value = -4407.5
$ws.Range("J21").Value = 7397.25
$ws.Range("K32").Value = 5318.3516
$ws.Range("H32").Value = 5022.7627
$ws.Range("I32").Value = 5318.3516
$ws.Range("M32").Value = -5031.3516
$ws.Range("M61").Value = -3774.4285
$ws.Range("H61").Value = 3986.4285
$ws.Range("I61").Value = 3986.4285
$ws.Range("K61").Value = 3986.4285
$ws.Range("H74").Value = 4254.8887
$ws.Range("L74").Value = 5286.143
$ws.Range("J74").Value = 5286.143
$ws.Range("N74").Value = -7034.143
$ws.Range("N77").Value = -35166.715
$ws.Range("H77").Value = 4254.8887
$ws.Range("J77").Value = 5286.143
$ws.Range("L77").Value = 26430.715
$ws.Range("M132").Value = -2131.7498
$ws.Range("K132").Value = 4661.7498
$ws.Range("H132").Value = 1505.5385
$ws.Range("I132").Value = 1553.9166
$ws.Range("K136").Value = 11959.2855
$ws.Range("H136").Value = 3986.4285
$ws.Range("M136").Value = -9409.2855
$ws.Range("I136").Value = 3986.4285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K99").Value = 3548
$ws.Range("I99").Value = 3548
$ws.Range("H99").Value = 4070.7778
$ws.Range("M99").Value = -2050
$ws.Range("H134").Value = 17809.727
$ws.Range("M134").Value = -9434.3334
$ws.Range("I134").Value = 3989.7778
$ws.Range("N134").Value = -245068.5
$ws.Range("L134").Value = 239998.5
$ws.Range("K134").Value = 11969.3334
$ws.Range("J134").Value = 79999.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("K99").Value = 2038187.6
$ws.Range("I99").Value = 2038187.6
$ws.Range("H99").Value = 14031203
$ws.Range("M99").Value = -2036689.6
$ws.Range("K122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("L122").Value = 13935
$ws.Range("N122").Value = -18835
$ws.Range("J122").Value = 4645
$ws.Range("H122").Value = 4645
$ws.Range("M122").ClearContents()
$ws.Range("K126").Value = 6114562.800000001
$ws.Range("H126").Value = 14031203
$ws.Range("I126").Value = 2038187.6
$ws.Range("M126").Value = -6112092.800000001
$ws.Range("N132").Value = -16010
$ws.Range("J132").Value = 3650
$ws.Range("L132").Value = 10950
$ws.Range("H132").Value = 18456.285
$ws.Range("H134").Value = 2354.516
$ws.Range("I134").Value = 2406.3333
$ws.Range("M134").Value = -4683.999899999999
$ws.Range("K134").Value = 7218.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M13").Value = -436.5
$ws.Range("K13").Value = 604.5
$ws.Range("I13").Value = 201.5
$ws.Range("H13").Value = 276

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K113").Value = 12148
$ws.Range("L113").Value = 16862.363
$ws.Range("H113").Value = 15198.471
$ws.Range("J113").Value = 16862.363
$ws.Range("I113").Value = 12148
$ws.Range("M113").Value = -9978
$ws.Range("N113").Value = -21202.363
$ws.Range("M132").Value = -23636.501
$ws.Range("K132").Value = 26166.501
$ws.Range("H132").Value = 9349.950000000001
$ws.Range("I132").Value = 8722.166999999999
$ws.Range("N141").Value = -55355
$ws.Range("J141").Value = 44995
$ws.Range("L141").Value = 44995
$ws.Range("H141").Value = 44995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K16").Value = 1724
$ws.Range("H16").Value = 1684
$ws.Range("M16").Value = -1554
$ws.Range("I16").Value = 1724
$ws.Range("H61").Value = 1500
$ws.Range("I61").Value = 1500
$ws.Range("K61").Value = 1500
$ws.Range("M61").Value = -1298
$ws.Range("H68").Value = 2536
$ws.Range("M68").Value = -1787
$ws.Range("I68").Value = 2536
$ws.Range("K68").Value = 2536
$ws.Range("I71").Value = 2536
$ws.Range("H71").Value = 2536
$ws.Range("M71").Value = -8936
$ws.Range("K71").Value = 12680
$ws.Range("K113").Value = 1500
$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1500
$ws.Range("M113").Value = 670
$ws.Range("K122").Value = 7274.7276
$ws.Range("I122").Value = 2424.9092
$ws.Range("M122").Value = -4824.7276
$ws.Range("H122").Value = 3905.5715
$ws.Range("M132").Value = -14511.155
$ws.Range("K132").Value = 17041.155
$ws.Range("H132").Value = 5454.8276
$ws.Range("I132").Value = 5680.385

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K113").Value = 2225.5713
$ws.Range("L113").Value = 1102.99998
$ws.Range("H113").Value = 629.6
$ws.Range("J113").Value = 367.66666
$ws.Range("I113").Value = 741.8570999999999
$ws.Range("M113").Value = -55.57129999999961
$ws.Range("N113").Value = -5442.999980000001
$ws.Range("K122").Value = 3863.4546
$ws.Range("I122").Value = 1287.8182
$ws.Range("M122").Value = -1413.4546
$ws.Range("L122").Value = 4156.2858
$ws.Range("N122").Value = -9056.2858
$ws.Range("H122").Value = 1325.7778
$ws.Range("J122").Value = 1385.4286
$ws.Range("M132").Value = -9062.706200000001
$ws.Range("K132").Value = 11592.7062
$ws.Range("H132").Value = 4093.2273
$ws.Range("I132").Value = 3864.2354
$ws.Range("K136").Value = 10627.5879
$ws.Range("H136").Value = 5057.5186
$ws.Range("J136").Value = 7633
$ws.Range("M136").Value = -8077.5879
$ws.Range("I136").Value = 3542.5293
$ws.Range("L136").Value = 22899
$ws.Range("N136").Value = -27999
$ws.Range("N138").Value = -85779.8
$ws.Range("H138").Value = 75499.8
$ws.Range("J138").Value = 75499.8
$ws.Range("L138").Value = 75499.8
